$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text change
$ws.Range("B1").Value = "Value (g)"

# Convert decimal values to text (shared strings), preserving default style
$ws.Range("B2").Formula = "=""30.5"""
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("B3").Formula = "=""26.76"""
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("B4").Formula = "=""51.3"""
$ws.Range("B4").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$ws.Range("B6").Formula = "=""43.97"""
$ws.Range("B6").Copy()
$ws.Range("B6").PasteSpecial(-4163)
$ws.Range("B7").Formula = "=""25.46"""
$ws.Range("B7").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("B8").Formula = "=""61.71"""
$ws.Range("B8").Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("B11").Formula = "=""5.88"""
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)
$ws.Range("B13").Formula = "=""49.64"""
$ws.Range("B13").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("B14").Formula = "=""25.46"""
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("B15").Formula = "=""37.6"""
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B16").Formula = "=""134.83"""
$ws.Range("B16").Copy()
$ws.Range("B16").PasteSpecial(-4163)
$ws.Range("B18").Formula = "=""83.99"""
$ws.Range("B18").Copy()
$ws.Range("B18").PasteSpecial(-4163)
$ws.Range("B20").Formula = "=""40.02"""
$ws.Range("B20").Copy()
$ws.Range("B20").PasteSpecial(-4163)
$ws.Range("B23").Formula = "=""13.2"""
$ws.Range("B23").Copy()
$ws.Range("B23").PasteSpecial(-4163)
$ws.Range("B24").Formula = "=""31.1"""
$ws.Range("B24").Copy()
$ws.Range("B24").PasteSpecial(-4163)
$ws.Range("B25").Formula = "=""97.93"""
$ws.Range("B25").Copy()
$ws.Range("B25").PasteSpecial(-4163)
$ws.Range("B26").Formula = "=""12.73"""
$ws.Range("B26").Copy()
$ws.Range("B26").PasteSpecial(-4163)
$ws.Range("B27").Formula = "=""61.51"""
$ws.Range("B27").Copy()
$ws.Range("B27").PasteSpecial(-4163)
$ws.Range("B29").Formula = "=""35.5"""
$ws.Range("B29").Copy()
$ws.Range("B29").PasteSpecial(-4163)
$ws.Range("B32").Formula = "=""50.97"""
$ws.Range("B32").Copy()
$ws.Range("B32").PasteSpecial(-4163)
$ws.Range("B33").Formula = "=""43.97"""
$ws.Range("B33").Copy()
$ws.Range("B33").PasteSpecial(-4163)
$ws.Range("B34").Formula = "=""18.08"""
$ws.Range("B34").Copy()
$ws.Range("B34").PasteSpecial(-4163)
$ws.Range("B35").Formula = "=""64.8"""
$ws.Range("B35").Copy()
$ws.Range("B35").PasteSpecial(-4163)
$ws.Range("B36").Formula = "=""2.8"""
$ws.Range("B36").Copy()
$ws.Range("B36").PasteSpecial(-4163)
$ws.Range("B37").Formula = "=""9.8"""
$ws.Range("B37").Copy()
$ws.Range("B37").PasteSpecial(-4163)
$ws.Range("B38").Formula = "=""21.3"""
$ws.Range("B38").Copy()
$ws.Range("B38").PasteSpecial(-4163)
$ws.Range("B39").Formula = "=""43.97"""
$ws.Range("B39").Copy()
$ws.Range("B39").PasteSpecial(-4163)

# Convert the three 3-decimal values to integer milli-units with #,##0 format
$ws.Range("B5").NumberFormat = "#,##0"
$ws.Range("B5").Value = 29385
$ws.Range("B17").NumberFormat = "#,##0"
$ws.Range("B17").Value = 16175
$ws.Range("B28").NumberFormat = "#,##0"
$ws.Range("B28").Value = 28175

# Column A width (bestFit)
$ws.Columns.Item(1).ColumnWidth = 17.6640625

# Selection change
$ws.Range("D7").Select()

$excel.CutCopyMode = $false
